$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '97.125.75'
$ws.Range("E2").Value = '  +0.18%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.711.63'
$ws.Range("E3").Value = '  +0.40%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.11'
$ws.Range("E5").Value = '  -2.58%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.90'
$ws.Range("E6").Value = '  +0.00%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '655.97'
$ws.Range("E7").Value = '  -1.94%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.431'
$ws.Range("E8").Value = '  -0.67%  '

# Row 9
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.07'
$ws.Range("E9").Value = '  -4.32%  '

# Row 10
$ws.Range("B10").Value = 'USDC'
$ws.Range("C10").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.999'
$ws.Range("E10").Value = '  +0.01%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.708.83'
$ws.Range("E11").Value = '  +0.42%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.43'
$ws.Range("E12").Value = '  -2.65%  '

# Row 13
$ws.Range("E13").Value = '  +1.17%  '

# Row 14
$ws.Range("E14").Value = '  +12.96%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.77'
$ws.Range("E15").Value = '  +2.61%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.402.45'
$ws.Range("E16").Value = '  +0.39%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.793.46'
$ws.Range("E17").Value = '  +0.09%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.94'
$ws.Range("E18").Value = '  -1.17%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.717.33'
$ws.Range("E19").Value = '  +0.55%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.25'
$ws.Range("E20").Value = '  +2.50%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.77'
$ws.Range("E21").Value = '  +1.09%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.509'
$ws.Range("E22").Value = '  -5.45%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '524.63'
$ws.Range("E23").Value = '  +1.67%  '

# Row 24
$ws.Range("E24").Value = '  -1.37%  '

# Row 25
$ws.Range("E25").Value = '  +1.38%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.96'
$ws.Range("E26").Value = '  +0.12%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '101.64'
$ws.Range("E27").Value = '  -0.03%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.196'
$ws.Range("E28").Value = '  +16.66%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '13.37'
$ws.Range("E29").Value = '  +1.71%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.22'
$ws.Range("E30").Value = '  +0.38%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.02'
$ws.Range("E31").Value = '  -1.86%  '

# Row 32
$ws.Range("E32").Value = '  +0.00%  '

# Row 33
$ws.Range("E33").Value = '  +1.67%  '

# Row 34
$ws.Range("E34").Value = '  +8.00%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.81%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '32.38'
$ws.Range("E36").Value = '  -2.10%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '649.33'
$ws.Range("E37").Value = '  +5.93%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.600'
$ws.Range("E38").Value = '  +1.69%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.86'
$ws.Range("E39").Value = '  +0.74%  '

# Row 40
$ws.Range("E40").Value = '  +0.03%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.84'
$ws.Range("E41").Value = '  +10.34%  '

# Row 42
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.80'
$ws.Range("E42").Value = '  -4.61%  '

# Row 43
$ws.Range("B43").Value = 'ImmutableX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.04'
$ws.Range("E43").Value = '  +4.48%  '

# Row 44
$ws.Range("E44").Value = '  -1.57%  '

# Row 45
$ws.Range("E45").Value = '  +0.00%  '

# Row 46
$ws.Range("E46").Value = '  +3.28%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0456'
$ws.Range("E47").Value = '  -1.55%  '

# Row 48
$ws.Range("E48").Value = '  -0.55%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.63'
$ws.Range("E49").Value = '  +0.03%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.56'
$ws.Range("E50").Value = '  -0.40%  '

# Row 51
$ws.Range("E51").Value = '  +0.59%  '
